$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The current header row (A1:D1) carries a bold/bordered/centered style that the
# new layout reuses for the new header cells (B1:E1) and for the new id column
# (A2:A5). Stash that style on a scratch cell far away so it survives while we
# clear and rebuild the A1:D4 block.
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# Wipe the old 4x4 block (values + formatting) completely; we rebuild it below.
$ws.Range("A1:D4").Clear()

# Re-apply the stashed header style to the new header row and id column.
$ws.Range("Z1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# ---- Header row (row 1): lang_code, code, name, is_active in B1:E1; A1 stays blank.
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "is_active"

# ---- Data rows 2-5
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eng"
$ws.Range("C2").Value = "FR"
$ws.Range("D2").Value = "Foreigner"
$ws.Range("E2").Value = $true

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "eng"
$ws.Range("C3").Value = "NFR"
$ws.Range("D3").Value = "Non-Foreigner"
$ws.Range("E3").Value = $true

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "fra"
$ws.Range("C4").Value = "FR"
$ws.Range("D4").Value = "Étranger"
$ws.Range("E4").Value = $true

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "fra"
$ws.Range("C5").Value = "NFR"
$ws.Range("D5").Value = "Non-étranger"
$ws.Range("E5").Value = $true

$ws.Range("A1").Select()
